# DEB_Fundamental_Indicators.xlsx edit
# - Rename the "indicator" header column to "feature" (A1 header cell on every
#   per-category sheet).
# - Add a new computed row to the "All" sheet.
# - Add a new "EBITAVG" indicator row to the "Income Statement" sheet.
# - Backfill the "priorImportance" (D) column on "Cash Flow Statement".
# - Restore view state (selections / active sheet) to match where the author
#   ended up.

$wb = $excel.ActiveWorkbook

$allSheet        = $wb.Worksheets.Item("All")
$solvency        = $wb.Worksheets.Item("Solvency and Fundamentals")
$riskPricing     = $wb.Worksheets.Item("Risk, Pricing and Valuation")
$profitability   = $wb.Worksheets.Item("Profitability and Management")
$incomeStatement = $wb.Worksheets.Item("Income Statement")
$cashFlow        = $wb.Worksheets.Item("Cash Flow Statement")
$balanceSheet    = $wb.Worksheets.Item("Balance Sheet")

# ---------------------------------------------------------------------------
# 1. Header rename: column A header "indicator" -> "feature" on every
#    per-category sheet (the "All" sheet is untouched - it doesn't carry
#    this header).
# ---------------------------------------------------------------------------
$solvency.Range("A1").Value        = "feature"
$riskPricing.Range("A1").Value     = "feature"
$profitability.Range("A1").Value   = "feature"
$incomeStatement.Range("A1").Value = "feature"
$cashFlow.Range("A1").Value        = "feature"
$balanceSheet.Range("A1").Value    = "feature"

# ---------------------------------------------------------------------------
# 2. "All" sheet: new row 199, a little scratch formula under the existing
#    indicator listing.
# ---------------------------------------------------------------------------
$allSheet.Range("A199").Formula = "=193 - 27"

# ---------------------------------------------------------------------------
# 3. "Income Statement": append a new indicator row (EBITAVG / 5 Year Average
#    EBIDTA / Annual / 0), matching the look of the existing data rows.
# ---------------------------------------------------------------------------
$incomeStatement.Range("A75").Value = "EBITAVG"
$incomeStatement.Range("B75").Value = "5 Year Average EBIDTA"
$incomeStatement.Range("C75").Value = "Annual"
$incomeStatement.Range("D75").Value = 0
$incomeStatement.Range("A75:C75").VerticalAlignment = -4108
$incomeStatement.Range("A75:C75").WrapText = $true

# ---------------------------------------------------------------------------
# 4. "Cash Flow Statement": backfill the priorImportance column (D) for all
#    existing data rows - it had been left blank before.
# ---------------------------------------------------------------------------
$cashFlowPriorImportance = @(1, 0, 0, 0, 0, 2, 2, 2, 2, 2)
for ($i = 0; $i -lt $cashFlowPriorImportance.Length; $i++) {
    $cashFlow.Cells.Item(2 + $i, 4).Value = $cashFlowPriorImportance[$i]
}

# ---------------------------------------------------------------------------
# 5. View state: selections on the sheets the author visited, finishing on
#    "Income Statement" so it ends up the active tab (matches activeTab going
#    from the old "Balance Sheet" tab to "Income Statement").
#    "Risk, Pricing and Valuation" and "Profitability and Management" are
#    left alone (no explicit selection either before or after).
# ---------------------------------------------------------------------------
$allSheet.Range("A116:C116").Select()
$solvency.Range("D24").Select()
$cashFlow.Range("D10").Select()
$incomeStatement.Range("D75").Select()
